$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 ("4. Include "Link Generator" as a feature"): Status Open -> Done, Assignee -> Arthur
$ws.Range("C4").Copy() | Out-Null
$ws.Range("C5").PasteSpecial(-4122) | Out-Null
$ws.Range("C5").Value = "Done"
$ws.Range("D5").Value = "Arthur"

# Row 7 ("6. Add feature to format urls"): Status Open -> Done, Assignee -> Arthur
$ws.Range("C4").Copy() | Out-Null
$ws.Range("C7").PasteSpecial(-4122) | Out-Null
$ws.Range("C7").Value = "Done"
$ws.Range("D7").Value = "Arthur"

# Row 21 ("20. Add review page"): Status In Progress -> Done
$ws.Range("C4").Copy() | Out-Null
$ws.Range("C21").PasteSpecial(-4122) | Out-Null
$ws.Range("C21").Value = "Done"

$excel.CutCopyMode = 0

# Move the active selection
$ws.Range("F14").Select() | Out-Null
